$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.760.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.49%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.549.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.84%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.89%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.11%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.24%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.96%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.95%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.939.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.523.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.02%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.40%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.776.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.39%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -9.67%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.92%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -8.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -10.60%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.56%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.32%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -10.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.50%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.081.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -10.29%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.95%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.793.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.81%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.37%  "
$ws.Range("E51").Style = "Normal"
